# Update "想去人数" (F column) values across the four sheets to reflect the
# newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5369
$ws1.Range("F6").Value = 5369
$ws1.Range("F11").Value = 1201
$ws1.Range("F12").Value = 6308
$ws1.Range("F14").Value = 80
$ws1.Range("F15").Value = 105
$ws1.Range("F16").Value = 3142
$ws1.Range("F18").Value = 106
$ws1.Range("F19").Value = 256
$ws1.Range("F20").Value = 4017
$ws1.Range("F24").Value = 3940
$ws1.Range("F25").Value = 191
$ws1.Range("F28").Value = 245
$ws1.Range("F29").Value = 257
$ws1.Range("F31").Value = 117
$ws1.Range("F32").Value = 124
$ws1.Range("F36").Value = 29
$ws1.Range("F37").Value = 6997
$ws1.Range("F38").Value = 32
$ws1.Range("F39").Value = 1153
$ws1.Range("F40").Value = 556
$ws1.Range("F43").Value = 1421
$ws1.Range("F46").Value = 2757
$ws1.Range("F47").Value = 322
$ws1.Range("F50").Value = 968

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 137
$ws2.Range("F12").Value = 4
$ws2.Range("F25").Value = 831

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 222

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 222
$ws4.Range("F8").Value = 5369
$ws4.Range("F9").Value = 5369
$ws4.Range("F15").Value = 1201
$ws4.Range("F16").Value = 6308
$ws4.Range("F18").Value = 80
$ws4.Range("F19").Value = 105
$ws4.Range("F20").Value = 3142
$ws4.Range("F22").Value = 106
$ws4.Range("F23").Value = 256
$ws4.Range("F24").Value = 4017
$ws4.Range("F25").Value = 3940
$ws4.Range("F26").Value = 191
$ws4.Range("F28").Value = 245
$ws4.Range("F29").Value = 257
$ws4.Range("F31").Value = 117
$ws4.Range("F32").Value = 124
$ws4.Range("F36").Value = 6997
$ws4.Range("F37").Value = 32
$ws4.Range("F38").Value = 1153
$ws4.Range("F39").Value = 556
$ws4.Range("F43").Value = 1421
$ws4.Range("F46").Value = 2758
$ws4.Range("F47").Value = 322
$ws4.Range("F49").Value = 968
